# Automatic update of files.
# The underlying data rows (3-7) were rotated: the record that used to sit
# in row 4 now sits in row 3, row 5's record moves to row 4, row 6's record
# moves to row 5, row 7's record moves to row 6, and row 3's original
# record wraps around into row 7. Apply the change by writing each row's
# new values directly (rather than shifting ranges), so columns that must
# stay blank/absent are explicitly cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3  (becomes old row 4: Spillkråka) ----
$ws.Range("A3").Value = 130937854
$ws.Range("B3").Value = 57881
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = "Spillkråka"
$ws.Range("G3").Value = "Dryocopus martius"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("J3").ClearContents()
$ws.Range("M3").Value = "färska spår"
$ws.Range("Q3").Value = 489668
$ws.Range("R3").Value = 7004128
$ws.Range("AC3").Value = "Rejäla hackspår, färska och äldre, I två levande granar och i ytlig grov rotdel."
$ws.Range("AF3").ClearContents()
$ws.Range("AJ3").Value = "gran"
$ws.Range("AK3").Value = "Picea abies"
$ws.Range("AM3").Value = "Trädstam på levande träd"
$ws.Range("AO3").Value = "Stem on living tree # Picea abies"

# ---- Row 4  (becomes old row 5: Tretåig hackspett, färska spår) ----
$ws.Range("A4").Value = 130937843
$ws.Range("B4").Value = 57884
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("Q4").Value = 489760
$ws.Range("R4").Value = 7004232
$ws.Range("AC4").Value = "Ringhack, färska och äldre, i riklig mängd längs flera meter högt upp på en granstam med spår av rikligt sav/kådaflöde."

# ---- Row 5  (becomes old row 6: Tretåig hackspett, äldre spår) ----
$ws.Range("A5").Value = 130937852
$ws.Range("M5").Value = "äldre spår"
$ws.Range("Q5").Value = 489520
$ws.Range("R5").Value = 7004161
$ws.Range("AC5").Value = "Ringhack, äldre, ytliga enstaka längs flera meter på en granstam vid kanten mot yngre skog."

# ---- Row 6  (becomes old row 7: Knärot) ----
$ws.Range("A6").Value = 130937863
$ws.Range("B6").Value = 99013
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("I6").Value = "'8"
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("K6").Value = "fullt utvecklade blad"
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 489799
$ws.Range("R6").Value = 7004245
$ws.Range("AC6").Value = "Minst 8 plantor inom ca 1 m2 yta. Grävdes varsamt fram under snötäcket. Det finns sannolikt betydligt mer knärot på fyndplatsen och i skogsbeståndet där fyndplatsen ligger."
$ws.Range("AF6").Value = "'"
$ws.Range("AH6").Value = "Barrskog"
$ws.Range("AJ6").ClearContents()
$ws.Range("AK6").ClearContents()
$ws.Range("AM6").ClearContents()
$ws.Range("AO6").ClearContents()

# ---- Row 7  (becomes old row 3: Revlummer) ----
$ws.Range("A7").Value = 130937857
$ws.Range("B7").Value = 97878
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 221945
$ws.Range("F7").Value = "Revlummer"
$ws.Range("G7").Value = "Lycopodium annotinum"
$ws.Range("H7").Value = "L."
$ws.Range("I7").Value = "'"
$ws.Range("J7").Value = "'"
$ws.Range("K7").Value = "'"
$ws.Range("Q7").Value = 489680
$ws.Range("R7").Value = 7004154
$ws.Range("AC7").ClearContents()
$ws.Range("AH7").Value = "Granskog"
